$d = $word.ActiveDocument

$target = "Avoir des modèles à suivre sur la toile"

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*$target*") {
        $rng = $p.Range

        # Pull this paragraph's exact canonical OOXML (preserves w14:paraId,
        # rsids, etc.) so we can surgically drop the w:sz/w:szCs overrides
        # without disturbing anything else about the paragraph.
        $full = $rng.WordOpenXML
        $bodyMatch = [System.Text.RegularExpressions.Regex]::Match(
            $full, '<w:body>(.*?)</w:body>', [System.Text.RegularExpressions.RegexOptions]::Singleline)
        $body = $bodyMatch.Groups[1].Value

        $pMatch = [System.Text.RegularExpressions.Regex]::Match(
            $body, '<w:p\b.*?</w:p>', [System.Text.RegularExpressions.RegexOptions]::Singleline)
        $parXml = $pMatch.Value

        # Remove the explicit (direct) run-size overrides so the text falls
        # back to the inherited/style size, both on the paragraph mark's
        # rPr (inside pPr) and on the run's own rPr.
        $parXml = $parXml -replace '<w:sz w:val="[0-9]+"/>', ''
        $parXml = $parXml -replace '<w:szCs w:val="[0-9]+"/>', ''

        $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
               '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData>' +
               '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:body>' + $parXml + '</w:body>' +
               '</w:document>' +
               '</pkg:xmlData></pkg:part></pkg:package>'

        $rng.InsertXML($xml)
        break
    }
}
